# repull data, push all data, mean calculation
# Update the dSF column (column F) values that changed after re-pulling
# and recomputing the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = -3
$ws.Cells.Item(4, 6).Value = -3
$ws.Cells.Item(5, 6).Value = -2
$ws.Cells.Item(7, 6).Value = -4
$ws.Cells.Item(12, 6).Value = -11
$ws.Cells.Item(13, 6).Value = -7
$ws.Cells.Item(15, 6).Value = 3
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(25, 6).Value = -6
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(29, 6).Value = 5
$ws.Cells.Item(33, 6).Value = -9
$ws.Cells.Item(34, 6).Value = 3
$ws.Cells.Item(35, 6).Value = -6
$ws.Cells.Item(36, 6).Value = -7
$ws.Cells.Item(37, 6).Value = -8
$ws.Cells.Item(38, 6).Value = -1
$ws.Cells.Item(39, 6).Value = -3
$ws.Cells.Item(40, 6).Value = -2
$ws.Cells.Item(42, 6).Value = -2
